$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16512555163073485"
$ws1.Range("B2").Value = "go_stims-16512555162663496.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555162903478.csv"
$ws1.Range("B4").Value = "go_stims-16512555162923493.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555163063507.csv"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-165125551903115"
$ws2.Range("B2").Value = "TB-16512555184361513.csv"
$ws2.Range("B3").Value = "OB-16512555169981527.csv"
$ws2.Range("B4").Value = "OB-16512555172831511.csv"
$ws2.Range("B5").Value = "TB-16512555190091507.csv"
$ws2.Range("B6").Value = "ZB-match_6-16512555167431502.csv"
$ws2.Range("B7").Value = "ZB-match_0-16512555168781517.csv"
$ws2.Range("B8").Value = "OB-16512555169661636.csv"
$ws2.Range("B9").Value = "ZB-match_2-1651255516720153.csv"
$ws2.Range("B10").Value = "TB-16512555184731512.csv"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16512555190331514"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16512555190791514"
$ws4.Range("B2").Value = "MM_stims-16512555190461593.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555190361552.csv"
$ws4.Range("B4").Value = "MM_stims-16512555190621538.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555190481532.csv"
$ws4.Range("B6").Value = "MM_stims-16512555190781536.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555190631526.csv"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16512555191571603"
$ws5.Range("B2").Value = "SAT_stims-1651255519086152.csv"
$ws5.Range("B3").Value = "SAT_stims-16512555191101534.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555191261523.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512555191421504.csv"
